# Updated cryptos list on Tue Oct 22 16:32:12 UTC 2024 with GitHub Actions
# Refreshes price (D) / 1h volume-change (E) columns for every coin row,
# and re-syncs three rows whose rank order shuffled (so name/link/price/
# change all moved together): rows 24-25 (Dai/Aptos) and rows 41-43/46-51
# (WhiteBITCoin/RenderToken/Stacks ... BabyDogeCoin/Aave/ARBITRUM/Filecoin/
# Optimism/Cronos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that round-trip as unambiguous numeric literals (e.g. a plain
# "594.59") must be forced to Text first, otherwise COM's auto-type-detection
# would silently convert them to a Double and drop the original formatting
# (trailing zeros, etc.) - exactly like typing into a cell already
# formatted as Text in the Excel UI.
$ws.Range("D2").Value = "66.997.98"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.617.73"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.59"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.73"
$ws.Range("E6").Value = "  +1.84%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").Value = "2.616.50"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.63"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").Value = "3.103.82"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000181"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").Value = "66.988.61"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "2.618.86"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.01"
$ws.Range("E19").Value = "  +3.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("E20").Value = "  +7.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.74"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("E22").Value = "  -0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.65"
$ws.Range("E23").Value = "  -2.38%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.33"
$ws.Range("E25").Value = "  +3.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.92"
$ws.Range("E26").Value = "  -4.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.54"
$ws.Range("E27").Value = "  -1.92%  "

$ws.Range("D28").Value = "2.756.24"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000100"
$ws.Range("E30").Value = "  -0.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "545.92"
$ws.Range("E31").Value = "  -0.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.89"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  -2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.136"
$ws.Range("E35").Value = "  +6.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  -3.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.11"
$ws.Range("E38").Value = "  +2.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.97"
$ws.Range("E39").Value = "  -2.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.366"
$ws.Range("E40").Value = "  -1.23%  "

$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.13"
$ws.Range("E41").Value = "  +1.39%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.19"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.80"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -3.43%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0297"
$ws.Range("E46").Value = "  +1.31%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.61"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.577"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.76"
$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0770"
$ws.Range("E51").Value = "  -0.76%  "
